$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.928091645240784
$ws.Range("B1").Value = 2.241205453872681
$ws.Range("C1").Value = 2.445765972137451
$ws.Range("D1").Value = 3.566776514053345
$ws.Range("E1").Value = 0.9629173874855042
